$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 with new category/title/weight columns
# (set E2 first so the shared-string table gets "PESEL" at index 19 and "cards" at index 20,
# matching the order new strings were authored in)
$ws.Range("E2").Value = "PESEL"
$ws.Range("D2").Value = "cards"
$ws.Range("H2").Value = 3
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 10
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 5

# Update the active selection to match the saved view state (cell F2 selected)
$ws.Range("F2").Select()
